# Update market / leve profit data refreshed by the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H5").Value = 95.84614999999999
$ws.Range("I5").Value = 96.7
$ws.Range("K5").Value = 96.7
$ws.Range("M5").Value = 18.3
$ws.Range("H6").Value = 1500.125
$ws.Range("I6").Value = 2100.2
$ws.Range("J6").Value = 500
$ws.Range("K6").Value = 6300.599999999999
$ws.Range("L6").Value = 1500
$ws.Range("M6").Value = -6188.599999999999
$ws.Range("N6").Value = -1724
$ws.Range("H18").Value = 1035
$ws.Range("I18").Value = 1035
$ws.Range("K18").Value = 1035
$ws.Range("M18").Value = -751
$ws.Range("H53").Value = 284.85715
$ws.Range("I53").Value = 89.09999999999999
$ws.Range("J53").Value = 462.81818
$ws.Range("K53").Value = 89.09999999999999
$ws.Range("L53").Value = 462.81818
$ws.Range("M53").Value = 547.9
$ws.Range("N53").Value = -1736.81818
$ws.Range("H58").Value = 2611.7
$ws.Range("I58").Value = 183.33333
$ws.Range("J58").Value = 6254.25
$ws.Range("K58").Value = 549.99999
$ws.Range("L58").Value = 18762.75
$ws.Range("M58").Value = -399.99999
$ws.Range("N58").Value = -19062.75
$ws.Range("H112").Value = 1259.2041
$ws.Range("I112").Value = 1020.125
$ws.Range("J112").Value = 1305.8536
$ws.Range("K112").Value = 3060.375
$ws.Range("L112").Value = 3917.5608
$ws.Range("M112").Value = -1952.375
$ws.Range("N112").Value = -6133.560799999999
$ws.Range("H137").Value = 1020.14703
$ws.Range("I137").Value = 743.76
$ws.Range("J137").Value = 1787.8889
$ws.Range("K137").Value = 2231.28
$ws.Range("L137").Value = 5363.6667
$ws.Range("M137").Value = 318.7200000000003
$ws.Range("N137").Value = -10463.6667

$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H2").Value = 1870
$ws.Range("I2").Value = 1400
$ws.Range("J2").Value = 2575
$ws.Range("K2").Value = 1400
$ws.Range("L2").Value = 2575
$ws.Range("M2").Value = -1287
$ws.Range("N2").Value = -2801
$ws.Range("H3").Value = 980
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 980
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 980
$ws.Range("M3").ClearContents()
$ws.Range("N3").Value = -1206
$ws.Range("H4").Value = 100
$ws.Range("I4").Value = 100
$ws.Range("K4").Value = 100
$ws.Range("M4").Value = 12
$ws.Range("H6").Value = 41303.637
$ws.Range("I6").Value = 80068
$ws.Range("J6").Value = 9000
$ws.Range("K6").Value = 80068
$ws.Range("L6").Value = 9000
$ws.Range("M6").Value = -79955
$ws.Range("N6").Value = -9226
$ws.Range("H7").Value = 145.93333
$ws.Range("I7").Value = 141.54546
$ws.Range("J7").Value = 158
$ws.Range("K7").Value = 141.54546
$ws.Range("L7").Value = 158
$ws.Range("M7").Value = -28.54545999999999
$ws.Range("N7").Value = -384
$ws.Range("H10").Value = 393.66666
$ws.Range("I10").Value = 440.5
$ws.Range("J10").Value = 300
$ws.Range("K10").Value = 440.5
$ws.Range("L10").Value = 300
$ws.Range("M10").Value = -301.5
$ws.Range("N10").Value = -578
$ws.Range("H11").Value = 52078
$ws.Range("I11").Value = 300
$ws.Range("J11").Value = 69337.336
$ws.Range("K11").Value = 300
$ws.Range("L11").Value = 69337.336
$ws.Range("M11").Value = -160
$ws.Range("N11").Value = -69617.336
$ws.Range("H12").Value = 195
$ws.Range("I12").Value = 195
$ws.Range("K12").Value = 195
$ws.Range("M12").Value = -25
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").ClearContents()
$ws.Range("H19").Value = 250
$ws.Range("I19").Value = 233.33333
$ws.Range("J19").Value = 300
$ws.Range("K19").Value = 233.33333
$ws.Range("L19").Value = 300
$ws.Range("M19").Value = -63.33332999999999
$ws.Range("N19").Value = -640
$ws.Range("H21").Value = 3900
$ws.Range("J21").Value = 3900
$ws.Range("L21").Value = 3900
$ws.Range("N21").Value = -4370
$ws.Range("H22").Value = 302.32
$ws.Range("I22").Value = 287.9
$ws.Range("J22").Value = 360
$ws.Range("K22").Value = 287.9
$ws.Range("L22").Value = 360
$ws.Range("M22").Value = 62.10000000000002
$ws.Range("N22").Value = -1060
$ws.Range("H23").Value = 70010
$ws.Range("J23").Value = 70010
$ws.Range("L23").Value = 70010
$ws.Range("N23").Value = -70490
$ws.Range("H24").Value = 250
$ws.Range("I24").Value = 233.33333
$ws.Range("J24").Value = 300
$ws.Range("K24").Value = 233.33333
$ws.Range("L24").Value = 300
$ws.Range("M24").Value = -63.33332999999999
$ws.Range("N24").Value = -640
$ws.Range("H27").Value = 70010
$ws.Range("J27").Value = 70010
$ws.Range("L27").Value = 70010
$ws.Range("N27").Value = -70394
$ws.Range("H31").Value = 1768.1111
$ws.Range("I31").Value = 1720.7307
$ws.Range("K31").Value = 1720.7307
$ws.Range("M31").Value = -1425.7307
$ws.Range("H34").Value = 1768.1111
$ws.Range("I34").Value = 1720.7307
$ws.Range("K34").Value = 1720.7307
$ws.Range("M34").Value = -1518.7307
$ws.Range("H86").Value = 3978.4285
$ws.Range("I86").Value = 2042.8667
$ws.Range("J86").Value = 6211.769
$ws.Range("K86").Value = 2042.8667
$ws.Range("L86").Value = 6211.769
$ws.Range("M86").Value = -919.8667
$ws.Range("N86").Value = -8457.769
$ws.Range("H89").Value = 3978.4285
$ws.Range("I89").Value = 2042.8667
$ws.Range("J89").Value = 6211.769
$ws.Range("K89").Value = 10214.3335
$ws.Range("L89").Value = 31058.845
$ws.Range("M89").Value = -4598.333500000001
$ws.Range("N89").Value = -42290.845

$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H17").Value = 383.33334
$ws.Range("J17").Value = 350
$ws.Range("L17").Value = 1050
$ws.Range("N17").Value = -1388
$ws.Range("H23").Value = 109.71429
$ws.Range("I23").Value = 300
$ws.Range("J23").Value = 78
$ws.Range("K23").Value = 900
$ws.Range("L23").Value = 234
$ws.Range("M23").Value = -665
$ws.Range("N23").Value = -704
$ws.Range("H64").Value = 2000
$ws.Range("I64").Value = 666.6667
$ws.Range("J64").Value = 3000
$ws.Range("K64").Value = 2000.0001
$ws.Range("L64").Value = 9000
$ws.Range("M64").Value = -1730.0001
$ws.Range("N64").Value = -9540
$ws.Range("H67").Value = 2000
$ws.Range("I67").Value = 666.6667
$ws.Range("J67").Value = 3000
$ws.Range("K67").Value = 2000.0001
$ws.Range("L67").Value = 9000
$ws.Range("M67").Value = -1064.0001
$ws.Range("N67").Value = -10872
$ws.Range("H131").Value = 782.78
$ws.Range("I131").Value = 407.14285
$ws.Range("J131").Value = 811.0538
$ws.Range("K131").Value = 1221.42855
$ws.Range("L131").Value = 2433.1614
$ws.Range("M131").Value = 3818.57145
$ws.Range("N131").Value = -12513.1614

$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H70").Value = 5300
$ws.Range("I70").Value = 4166.6665
$ws.Range("J70").Value = 7000
$ws.Range("K70").Value = 4166.6665
$ws.Range("L70").Value = 7000
$ws.Range("M70").Value = -3896.6665
$ws.Range("N70").Value = -7540
$ws.Range("H73").Value = 5300
$ws.Range("I73").Value = 4166.6665
$ws.Range("J73").Value = 7000
$ws.Range("K73").Value = 4166.6665
$ws.Range("M73").Value = -3230.6665
$ws.Range("N73").Value = -8872

